$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the obsolete "ci" guide row (row 4). Every row below shifts up by
#    one, which is exactly what the target workbook looks like (one fewer
#    row overall, "ci" no longer present anywhere).
$ws.Rows(4).Delete() | Out-Null

# 2. New "Generated Meta Tag" column (F) header, mirroring the existing
#    "Generated Title Tag" column (E).
$ws.Range("F1").Value = "Generated Meta Tag"

# 3. New "Topic" values (column C) for the first few guides, plus a test
#    description (column D) on the Accessibility row.
$ws.Range("D2").Value = "Test description"
$ws.Range("C2").Value = "Learn Accessibility"
$ws.Range("C3").Value = "Learn Browser, Website, and Web App Performance"
$ws.Range("C4").Value = "Learn Code Editors"
$ws.Range("C5").Value = "Learn the Command Line Interface"
$ws.Range("C6").Value = "Learn CSS Fonts & Icons"
$ws.Range("C7").Value = "Learn CSS Fundamentals"
$ws.Range("C8").Value = "Learn CSS in JavaScript"
$ws.Range("C9").Value = "Learn CSS Layout"

$ws.Range("E2").Formula = '="<title>Learning "&C2&" Resources - Front-End Developer Learning Roadmap</title>"'
$ws.Range("F2").Formula = "=" + '"<meta name=' + "'description' content=''" + '"&D2&"' + "'>" + '"'

$ws.Range("E3:E59").Formula = '="<title>Learning "&C3&" Resources - Front-End Developer Learning Roadmap</title>"'
$ws.Range("F3:F59").Formula = "=" + '"<meta name=' + "'description' content=''" + '"&D3&"' + "'>" + '"'

# 4. Column widths (characters) and the final active-cell selection.
$ws.Columns("C").ColumnWidth = 54.33
$ws.Columns("E").ColumnWidth = 104
$ws.Columns("F").ColumnWidth = 84.67

$ws.Range("C9").Select() | Out-Null
